$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename shared-string header labels: _old -> _FV2310, _new -> _FV2404 ---
for ($col = 1; $col -le 21; $col++) {
    $cell = $ws.Cells.Item(1, $col)
    $v = $cell.Value()
    if ($v -ne $null) {
        if ($v.EndsWith("_old")) {
            $cell.Value = $v.Substring(0, $v.Length - 4) + "_FV2310"
        } elseif ($v.EndsWith("_new")) {
            $cell.Value = $v.Substring(0, $v.Length - 4) + "_FV2404"
        }
    }
}

# --- 2. Freeze the header row ---
$ws.Activate()
$excel.ActiveWindow.FreezePanes = $false
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

# --- 3. Turn the used range into a table (ListObject) ---
$rng = $ws.Range("A1:U89")
$lo = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $rng, $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$lo.Name = "Table1"
